# exampleInterview.xlsx: add support for a more powerful tag format.
#
# - sheet "links" is renamed "tags" and becomes a lookup table of tag
#   groups/values (group/name/label/imagePath) instead of simple from/to
#   links.
# - sheet "interview" gains a "type" column (begin branch / end branch /
#   Follow-on question) and a "tags" column, replacing the old
#   from/to/noNextLink linking model, and a couple of copy fixes.
# - the "tags" sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

$interview = $wb.Worksheets.Item("interview")
$links     = $wb.Worksheets.Item("links")

# ------------------------------------------------------------------
# Rename "links" -> "tags"
# ------------------------------------------------------------------
$links.Name = "tags"
$tags = $links

# ------------------------------------------------------------------
# "interview" sheet: wipe the old from/to table and rebuild it as the
# new type/name/label/tags table.
# ------------------------------------------------------------------
$interview.Cells.Clear()

$interview.Columns.Item(1).ColumnWidth = 13.5
$interview.Columns.Item(2).ColumnWidth = 13.5
$interview.Columns.Item(3).ColumnWidth = 41.0

$interview.Range("A1").Value = "type"
$interview.Range("B1").Value = "name"
$interview.Range("C1").Value = "label"
$interview.Range("D1").Value = "tags"
$interview.Rows.Item(1).AutoFit()

$interview.Range("B2").Value = "start"
$interview.Range("C2").Value = "Tell me about some of the tasks you use computers for?"
$interview.Rows.Item(2).RowHeight = 32.25

$interview.Range("A3").Value = "begin branch"
$interview.Range("C3").Value = "Follow-on question"
$interview.Rows.Item(3).RowHeight = 38.25

$interview.Range("C4").Value = "What tasks do you find the most challenging or time consuming?"
$interview.Range("D4").Value = "task_tags"
$interview.Rows.Item(4).RowHeight = 38.25

$interview.Range("A6").Value = "end branch"
$interview.Rows.Item(6).RowHeight = 38.25

$interview.Range("A7").Value = "begin branch"
$interview.Range("C7").Value = "Subject does not do any computer related tasks."
$interview.Rows.Item(7).RowHeight = 32.25

$interview.Range("A8").Value = "end branch"
$interview.Rows.Item(8).RowHeight = 32.25

$interview.Range("C9").Value = "What was your first experience with a computer?"
$interview.Range("D9").Value = "none"
$interview.Rows.Item(9).RowHeight = 12.75

$interview.Range("D29").Select()

# ------------------------------------------------------------------
# "tags" sheet: wipe the old from/to/tags table and rebuild it as the
# new group/name/label/imagePath tag catalogue.
# ------------------------------------------------------------------
$tags.Cells.Clear()

$tags.Columns.Item(2).ColumnWidth = 14.166666666666666
$tags.Columns.Item(3).ColumnWidth = 12.666666666666666
$tags.Columns.Item(4).ColumnWidth = 10.0

$tags.Range("A1").Value = "group"
$tags.Range("B1").Value = "name"
$tags.Range("C1").Value = "label"
$tags.Range("D1").Value = "imagePath"
$tags.Rows.Item(1).AutoFit()

$tags.Range("A2").Value = "default"
$tags.Range("B2").Value = "important"
$tags.Range("C2").Value = "important"
$tags.Rows.Item(2).AutoFit()

$tags.Range("A3").Value = "default"
$tags.Range("B3").Value = "tangential"
$tags.Range("C3").Value = "tangential"
$tags.Rows.Item(3).AutoFit()

$tags.Range("A4").Value = "default"
$tags.Range("B4").Value = "star"
$tags.Range("D4").Value = "star.png"
$tags.Rows.Item(4).AutoFit()

$tags.Range("A5").Value = "task_tags"
$tags.Range("B5").Value = "programming"
$tags.Range("C5").Value = "programming"
$tags.Rows.Item(5).RowHeight = 25.5

$tags.Range("A6").Value = "task_tags"
$tags.Range("B6").Value = "data entry"
$tags.Range("C6").Value = "data entry"
$tags.Rows.Item(6).AutoFit()

$tags.Range("C10").Select()
$tags.Activate()
